$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 1024.75
$ws.Range("I2").Value = 379.4
$ws.Range("K2").Value = 379.4
$ws.Range("M2").Value = -266.4

$ws.Range("H12").Value = 145.71428
$ws.Range("I12").Value = 136.66667
$ws.Range("J12").Value = 200
$ws.Range("K12").Value = 136.66667
$ws.Range("L12").Value = 200
$ws.Range("M12").Value = 33.33332999999999
$ws.Range("N12").Value = -540

$ws.Range("H70").Value = 3250
$ws.Range("I70").Value = 3250
$ws.Range("J70").Value = 0
$ws.Range("K70").Value = 9750
$ws.Range("L70").Value = 0
$ws.Range("M70").Value = -9480
$ws.Range("N70").ClearContents()

$ws.Range("H73").Value = 3250
$ws.Range("I73").Value = 3250
$ws.Range("J73").Value = 0
$ws.Range("K73").Value = 9750
$ws.Range("L73").Value = 0
$ws.Range("M73").Value = -8814
$ws.Range("N73").ClearContents()

$ws.Range("H74").Value = 7311.385
$ws.Range("I74").Value = 3504.8
$ws.Range("K74").Value = 3504.8
$ws.Range("M74").Value = -2568.8

$ws.Range("H77").Value = 7311.385
$ws.Range("I77").Value = 3504.8
$ws.Range("K77").Value = 17524
$ws.Range("M77").Value = -12844

$ws.Range("H86").Value = 4499.1665
$ws.Range("J86").Value = 5248.75
$ws.Range("L86").Value = 5248.75
$ws.Range("N86").Value = -7494.75

$ws.Range("H89").Value = 4499.1665
$ws.Range("J89").Value = 5248.75
$ws.Range("L89").Value = 26243.75
$ws.Range("N89").Value = -37475.75

$ws.Range("H92").Value = 2203.818
$ws.Range("I92").Value = 2203.818
$ws.Range("J92").Value = 0
$ws.Range("K92").Value = 2203.818
$ws.Range("L92").Value = 0
$ws.Range("M92").Value = -955.8180000000002
$ws.Range("N92").ClearContents()

$ws.Range("H107").Value = 789.48486
$ws.Range("I107").Value = 769.0741
$ws.Range("J107").Value = 881.3333
$ws.Range("K107").Value = 769.0741
$ws.Range("L107").Value = 881.3333
$ws.Range("M107").Value = 1150.9259
$ws.Range("N107").Value = -4721.3333

$ws.Range("H137").Value = 3778.9546
$ws.Range("I137").Value = 4603.5454
$ws.Range("J137").Value = 2954.3635
$ws.Range("K137").Value = 13810.6362
$ws.Range("L137").Value = 8863.0905
$ws.Range("M137").Value = -11260.6362
$ws.Range("N137").Value = -13963.0905

$ws.Range("H138").Value = 5244.7856
$ws.Range("J138").Value = 5340.4614
$ws.Range("L138").Value = 16021.3842
$ws.Range("N138").Value = -26301.3842

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1289.7693
$ws.Range("I2").Value = 858.375
$ws.Range("J2").Value = 1980
$ws.Range("K2").Value = 858.375
$ws.Range("L2").Value = 1980
$ws.Range("M2").Value = -745.375
$ws.Range("N2").Value = -2206

$ws.Range("J4").Value = 400
$ws.Range("L4").Value = 400
$ws.Range("N4").Value = -632

$ws.Range("H5").Value = 500500
$ws.Range("I5").Value = 500500
$ws.Range("K5").Value = 500500
$ws.Range("M5").Value = -500388

$ws.Range("H31").Value = 7364.923
$ws.Range("I31").Value = 5478.8335
$ws.Range("J31").Value = 29998
$ws.Range("K31").Value = 5478.8335
$ws.Range("L31").Value = 29998
$ws.Range("M31").Value = -5184.8335
$ws.Range("N31").Value = -30586

$ws.Range("H32").Value = 12611.07
$ws.Range("I32").Value = 9260.641
$ws.Range("J32").Value = 39719.09
$ws.Range("K32").Value = 9260.641
$ws.Range("L32").Value = 39719.09
$ws.Range("M32").Value = -8973.641
$ws.Range("N32").Value = -40293.09

$ws.Range("H37").Value = 16955.2

$ws.Range("H45").Value = 6455.3
$ws.Range("I45").Value = 9654.917
$ws.Range("J45").Value = 1655.875
$ws.Range("K45").Value = 9654.917
$ws.Range("L45").Value = 1655.875
$ws.Range("M45").Value = -9277.917
$ws.Range("N45").Value = -2409.875

$ws.Range("H46").Value = 13114.286
$ws.Range("J46").Value = 13800.167
$ws.Range("L46").Value = 13800.167
$ws.Range("N46").Value = -14438.167

$ws.Range("H55").Value = 2000
$ws.Range("J55").Value = 0
$ws.Range("L55").Value = 0
$ws.Range("N55").ClearContents()

$ws.Range("H61").Value = 449529.5
$ws.Range("I61").Value = 4458.0835
$ws.Range("J61").Value = 1339672.4
$ws.Range("K61").Value = 4458.0835
$ws.Range("L61").Value = 1339672.4
$ws.Range("M61").Value = -4246.0835
$ws.Range("N61").Value = -1340096.4

$ws.Range("H70").Value = 0
$ws.Range("J70").Value = 0
$ws.Range("L70").Value = 0
$ws.Range("N70").ClearContents()

$ws.Range("H73").Value = 0
$ws.Range("J73").Value = 0
$ws.Range("L73").Value = 0
$ws.Range("N73").ClearContents()

$ws.Range("H74").Value = 80866
$ws.Range("I74").Value = 102600.1
$ws.Range("K74").Value = 102600.1
$ws.Range("M74").Value = -101726.1

$ws.Range("H77").Value = 80866
$ws.Range("I77").Value = 102600.1
$ws.Range("K77").Value = 513000.5
$ws.Range("M77").Value = -508632.5

$ws.Range("H102").Value = 1429.9048
$ws.Range("I102").Value = 1479.3889
$ws.Range("J102").Value = 1133
$ws.Range("K102").Value = 1479.3889
$ws.Range("L102").Value = 1133
$ws.Range("M102").Value = 142.6111000000001
$ws.Range("N102").Value = -4377

$ws.Range("H116").Value = 1289.7693
$ws.Range("I116").Value = 858.375
$ws.Range("J116").Value = 1980
$ws.Range("K116").Value = 858.375
$ws.Range("L116").Value = 1980
$ws.Range("M116").Value = 1435.625
$ws.Range("N116").Value = -6568

$ws.Range("H132").Value = 3959.9375
$ws.Range("I132").Value = 3613.4333
$ws.Range("K132").Value = 10840.2999
$ws.Range("M132").Value = -8310.2999

$ws.Range("H136").Value = 449529.5
$ws.Range("I136").Value = 4458.0835
$ws.Range("J136").Value = 1339672.4
$ws.Range("K136").Value = 13374.2505
$ws.Range("L136").Value = 4019017.2
$ws.Range("M136").Value = -10824.2505
$ws.Range("N136").Value = -4024117.2

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1289.7693
$ws.Range("I3").Value = 858.375
$ws.Range("J3").Value = 1980
$ws.Range("K3").Value = 858.375
$ws.Range("L3").Value = 1980
$ws.Range("M3").Value = -744.375
$ws.Range("N3").Value = -2208

$ws.Range("H4").Value = 500500
$ws.Range("I4").Value = 500500
$ws.Range("K4").Value = 500500
$ws.Range("M4").Value = -500385

$ws.Range("H22").Value = 2828.4285
$ws.Range("I22").Value = 949.75
$ws.Range("K22").Value = 949.75
$ws.Range("M22").Value = -776.75

$ws.Range("H29").Value = 10533
$ws.Range("I29").Value = 10954.182
$ws.Range("K29").Value = 10954.182
$ws.Range("M29").Value = -10665.182

$ws.Range("H44").Value = 29998.25
$ws.Range("J44").Value = 29998.25
$ws.Range("L44").Value = 29998.25
$ws.Range("N44").Value = -30992.25

$ws.Range("H86").Value = 2171.7273
$ws.Range("I86").Value = 2177.111
$ws.Range("J86").Value = 2147.5
$ws.Range("K86").Value = 2177.111
$ws.Range("L86").Value = 2147.5
$ws.Range("M86").Value = -1054.111
$ws.Range("N86").Value = -4393.5

$ws.Range("H89").Value = 2171.7273
$ws.Range("I89").Value = 2177.111
$ws.Range("J89").Value = 2147.5
$ws.Range("K89").Value = 10885.555
$ws.Range("L89").Value = 10737.5
$ws.Range("M89").Value = -5269.555
$ws.Range("N89").Value = -21969.5

$ws.Range("H94").Value = 553.7143
$ws.Range("I94").Value = 589
$ws.Range("J94").Value = 403.75
$ws.Range("K94").Value = 589
$ws.Range("L94").Value = 403.75
$ws.Range("M94").Value = -138
$ws.Range("N94").Value = -1305.75

$ws.Range("H99").Value = 8899.4
$ws.Range("I99").Value = 1498.5
$ws.Range("J99").Value = 13833.333
$ws.Range("K99").Value = 1498.5
$ws.Range("L99").Value = 13833.333
$ws.Range("M99").Value = -0.5
$ws.Range("N99").Value = -16829.333

$ws.Range("H105").Value = 1406.125
$ws.Range("I105").Value = 1406.125
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 1406.125
$ws.Range("L105").Value = 0
$ws.Range("M105").Value = 340.875
$ws.Range("N105").ClearContents()

$ws.Range("H107").Value = 5000
$ws.Range("I107").Value = 5000
$ws.Range("K107").Value = 5000
$ws.Range("M107").Value = -3080

$ws.Range("H132").Value = 100779.5
$ws.Range("J132").Value = 100779.5
$ws.Range("L132").Value = 100779.5
$ws.Range("N132").Value = -110899.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 254.07143
$ws.Range("I7").Value = 40.88889
$ws.Range("J7").Value = 637.8
$ws.Range("K7").Value = 40.88889
$ws.Range("L7").Value = 637.8
$ws.Range("M7").Value = 72.11111
$ws.Range("N7").Value = -863.8

$ws.Range("H31").Value = 2803.9512
$ws.Range("I31").Value = 2425.8918
$ws.Range("J31").Value = 6301
$ws.Range("K31").Value = 2425.8918
$ws.Range("L31").Value = 6301
$ws.Range("M31").Value = -2130.8918
$ws.Range("N31").Value = -6891

$ws.Range("H34").Value = 2803.9512
$ws.Range("I34").Value = 2425.8918
$ws.Range("J34").Value = 6301
$ws.Range("K34").Value = 2425.8918
$ws.Range("L34").Value = 6301
$ws.Range("M34").Value = -2223.8918
$ws.Range("N34").Value = -6705

$ws.Range("H62").Value = 5046.875
$ws.Range("I62").Value = 3625
$ws.Range("K62").Value = 3625
$ws.Range("M62").Value = -3001

$ws.Range("H65").Value = 5046.875
$ws.Range("I65").Value = 3625
$ws.Range("K65").Value = 18125
$ws.Range("M65").Value = -15005

$ws.Range("H69").Value = 18727
$ws.Range("I69").Value = 3090.5
$ws.Range("K69").Value = 3090.5
$ws.Range("M69").Value = -2341.5

$ws.Range("H72").Value = 18727
$ws.Range("I72").Value = 3090.5
$ws.Range("K72").Value = 9271.5
$ws.Range("M72").Value = -5527.5

$ws.Range("H86").Value = 3261.923
$ws.Range("I86").Value = 2913.375
$ws.Range("K86").Value = 2913.375
$ws.Range("M86").Value = -1790.375

$ws.Range("H89").Value = 3261.923
$ws.Range("I89").Value = 2913.375
$ws.Range("K89").Value = 14566.875
$ws.Range("M89").Value = -8950.875

$ws.Range("H105").Value = 845
$ws.Range("I105").Value = 845
$ws.Range("K105").Value = 845
$ws.Range("M105").Value = 902

$ws.Range("H107").Value = 778.1
$ws.Range("I107").Value = 673.9091
$ws.Range("K107").Value = 673.9091
$ws.Range("M107").Value = 1246.0909

$ws.Range("H124").Value = 39999.668
$ws.Range("J124").Value = 39999.668
$ws.Range("L124").Value = 39999.668
$ws.Range("N124").Value = -44909.668

$ws.Range("H132").Value = 2545.7778
$ws.Range("I132").Value = 2166.5715
$ws.Range("K132").Value = 6499.7145
$ws.Range("M132").Value = -3969.7145

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 100
$ws.Range("I2").Value = 75
$ws.Range("J2").Value = 125
$ws.Range("K2").Value = 450
$ws.Range("L2").Value = 750
$ws.Range("M2").Value = -337
$ws.Range("N2").Value = -976

$ws.Range("H33").Value = 153.96
$ws.Range("I33").Value = 126.82353
$ws.Range("J33").Value = 211.625
$ws.Range("K33").Value = 760.94118
$ws.Range("L33").Value = 1269.75
$ws.Range("M33").Value = -477.94118
$ws.Range("N33").Value = -1835.75

$ws.Range("H38").Value = 57.05
$ws.Range("I38").Value = 52.384617
$ws.Range("K38").Value = 157.153851
$ws.Range("M38").Value = 189.846149

$ws.Range("H105").Value = 5000
$ws.Range("I105").Value = 5000
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 15000
$ws.Range("L105").Value = 0
$ws.Range("N105").ClearContents()
$ws.Range("M105").Value = -12379

$ws.Range("H113").Value = 5954281
$ws.Range("I113").Value = 11905962
$ws.Range("J113").Value = 2600
$ws.Range("K113").Value = 35717886
$ws.Range("L113").Value = 7800
$ws.Range("M113").Value = -35715716
$ws.Range("N113").Value = -12140

$ws.Range("H129").Value = 11000696
$ws.Range("I129").Value = 14143357
$ws.Range("J129").Value = 1383
$ws.Range("K129").Value = 42430071
$ws.Range("L129").Value = 4149
$ws.Range("M129").Value = -42425071
$ws.Range("N129").Value = -14149

$ws.Range("H131").Value = 4215.0757
$ws.Range("I131").Value = 1279.8
$ws.Range("J131").Value = 4520.8335
$ws.Range("K131").Value = 3839.4
$ws.Range("L131").Value = 13562.5005
$ws.Range("M131").Value = 1200.6
$ws.Range("N131").Value = -23642.5005

$ws.Range("H134").Value = 9165.833
$ws.Range("I134").Value = 5998.5713
$ws.Range("K134").Value = 17995.7139
$ws.Range("M134").Value = -12925.7139

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H39").Value = 0
$ws.Range("J39").Value = 0
$ws.Range("L39").Value = 0
$ws.Range("N39").ClearContents()

$ws.Range("H49").Value = 26718.143
$ws.Range("I49").Value = 25027
$ws.Range("K49").Value = 25027
$ws.Range("M49").Value = -24843

$ws.Range("H63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("N63").ClearContents()

$ws.Range("H66").Value = 0
$ws.Range("J66").Value = 0
$ws.Range("L66").Value = 0
$ws.Range("N66").ClearContents()

$ws.Range("H80").Value = 3420
$ws.Range("I80").Value = 2208
$ws.Range("K80").Value = 2208
$ws.Range("M80").Value = -1210

$ws.Range("H83").Value = 3420
$ws.Range("I83").Value = 2208
$ws.Range("K83").Value = 11040
$ws.Range("M83").Value = -6048

$ws.Range("H97").Value = 2218.7693
$ws.Range("I97").Value = 2342.6667
$ws.Range("K97").Value = 2342.6667
$ws.Range("M97").Value = -1846.6667

$ws.Range("H102").Value = 43479696
$ws.Range("I102").Value = 1759.8667
$ws.Range("K102").Value = 1759.8667
$ws.Range("M102").Value = -137.8667

$ws.Range("H122").Value = 2580.6206
$ws.Range("I122").Value = 2493.7144
$ws.Range("K122").Value = 7481.1432
$ws.Range("M122").Value = -5031.1432

$ws.Range("H132").Value = 3189.7144
$ws.Range("I132").Value = 2191.889
$ws.Range("J132").Value = 4985.8
$ws.Range("K132").Value = 6575.667
$ws.Range("L132").Value = 14957.4
$ws.Range("M132").Value = -4045.667
$ws.Range("N132").Value = -20017.4

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6300.222
$ws.Range("I7").Value = 5956.7144
$ws.Range("K7").Value = 5956.7144
$ws.Range("M7").Value = -5844.7144

$ws.Range("H22").Value = 4511.1763
$ws.Range("I22").Value = 3781.6667
$ws.Range("J22").Value = 4909.091
$ws.Range("K22").Value = 3781.6667
$ws.Range("L22").Value = 4909.091
$ws.Range("M22").Value = -3486.6667
$ws.Range("N22").Value = -5499.091

$ws.Range("H27").Value = 4511.1763
$ws.Range("I27").Value = 3781.6667
$ws.Range("J27").Value = 4909.091
$ws.Range("K27").Value = 3781.6667
$ws.Range("L27").Value = 4909.091
$ws.Range("M27").Value = -3674.6667
$ws.Range("N27").Value = -5123.091

$ws.Range("H55").Value = 194.16667
$ws.Range("I55").Value = 93.2
$ws.Range("J55").Value = 699
$ws.Range("K55").Value = 93.2
$ws.Range("L55").Value = 699
$ws.Range("M55").Value = 79.8
$ws.Range("N55").Value = -1045

$ws.Range("H126").Value = 6300.222
$ws.Range("I126").Value = 5956.7144
$ws.Range("K126").Value = 17870.1432
$ws.Range("M126").Value = -15400.1432

$ws.Range("H132").Value = 4408.579
$ws.Range("I132").Value = 4234.077
$ws.Range("J132").Value = 4786.6665
$ws.Range("K132").Value = 12702.231
$ws.Range("L132").Value = 14359.9995
$ws.Range("M132").Value = -10172.231
$ws.Range("N132").Value = -19419.9995

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 103333.336
$ws.Range("J46").Value = 103333.336
$ws.Range("L46").Value = 103333.336
$ws.Range("N46").Value = -103795.336

$ws.Range("H52").Value = 23359.4
$ws.Range("I52").Value = 23359.4
$ws.Range("K52").Value = 23359.4
$ws.Range("M52").Value = -23133.4

$ws.Range("H62").Value = 26499.8
$ws.Range("I62").Value = 24833
$ws.Range("K62").Value = 24833
$ws.Range("M62").Value = -24209

$ws.Range("H65").Value = 26499.8
$ws.Range("I65").Value = 24833
$ws.Range("K65").Value = 124165
$ws.Range("M65").Value = -121045

$ws.Range("H107").Value = 1030.25
$ws.Range("I107").Value = 795.06665
$ws.Range("J107").Value = 1735.8
$ws.Range("K107").Value = 2385.19995
$ws.Range("L107").Value = 5207.4
$ws.Range("M107").Value = -465.1999500000002
$ws.Range("N107").Value = -9047.4

$ws.Range("H126").Value = 10895.786
$ws.Range("I126").Value = 13495.1
$ws.Range("J126").Value = 4397.5
$ws.Range("K126").Value = 40485.3
$ws.Range("L126").Value = 13192.5
$ws.Range("M126").Value = -38015.3
$ws.Range("N126").Value = -18132.5

$ws.Range("H132").Value = 4415.0386
$ws.Range("I132").Value = 4044.0557
$ws.Range("J132").Value = 5249.75
$ws.Range("K132").Value = 12132.1671
$ws.Range("L132").Value = 15749.25
$ws.Range("M132").Value = -9602.167099999999
$ws.Range("N132").Value = -20809.25

$ws.Range("H134").Value = 103333.336
$ws.Range("J134").Value = 103333.336
$ws.Range("L134").Value = 310000.008
$ws.Range("N134").Value = -315070.008

$ws.Range("H136").Value = 2491.8076
$ws.Range("I136").Value = 2190
$ws.Range("J136").Value = 3170.875
$ws.Range("K136").Value = 6570
$ws.Range("L136").Value = 9512.625
$ws.Range("M136").Value = -4020
$ws.Range("N136").Value = -14612.625
